# Auto-generated edit script applying Ifrit_Profits.xlsx market-data refresh
# Updates currentAveragePrice* / LevePrice* / LeveProfit* columns (H:N) for specific leve rows
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 138
$ws.Cells.Item(138, 8).Value = 2535.348
$ws.Cells.Item(138, 9).Value = 0
$ws.Cells.Item(138, 10).Value = 2535.348
$ws.Cells.Item(138, 11).Value = 0
$ws.Cells.Item(138, 12).ClearContents()
$ws.Cells.Item(138, 13).Value = 7606.044
$ws.Cells.Item(138, 14).Value = -17886.044

$ws = $wb.Worksheets.Item("ARM")
# ARM row 32
$ws.Cells.Item(32, 8).Value = 5043.7847
$ws.Cells.Item(32, 9).Value = 5334.654
$ws.Cells.Item(32, 11).Value = 5334.654
$ws.Cells.Item(32, 13).Value = -5047.654

# ARM row 45
$ws.Cells.Item(45, 8).Value = 823.8
$ws.Cells.Item(45, 9).Value = 790
$ws.Cells.Item(45, 10).Value = 846.3333
$ws.Cells.Item(45, 11).Value = 790
$ws.Cells.Item(45, 12).Value = 846.3333
$ws.Cells.Item(45, 13).Value = -413
$ws.Cells.Item(45, 14).Value = -1600.3333

# ARM row 74
$ws.Cells.Item(74, 8).Value = 4586.394
$ws.Cells.Item(74, 9).Value = 1165.4615
$ws.Cells.Item(74, 11).Value = 1165.4615
$ws.Cells.Item(74, 13).Value = -291.4614999999999

# ARM row 77
$ws.Cells.Item(77, 8).Value = 4586.394
$ws.Cells.Item(77, 9).Value = 1165.4615
$ws.Cells.Item(77, 11).Value = 5827.307499999999
$ws.Cells.Item(77, 13).Value = -1459.307499999999

$ws = $wb.Worksheets.Item("BSM")
# BSM row 99
$ws.Cells.Item(99, 8).Value = 1726.8
$ws.Cells.Item(99, 9).Value = 1453
$ws.Cells.Item(99, 10).Value = 2137.5
$ws.Cells.Item(99, 11).Value = 1453
$ws.Cells.Item(99, 12).Value = 2137.5
$ws.Cells.Item(99, 13).Value = 45
$ws.Cells.Item(99, 14).Value = -5133.5

$ws = $wb.Worksheets.Item("CRP")
# CRP row 16
$ws.Cells.Item(16, 8).Value = 1502.625
$ws.Cells.Item(16, 9).Value = 1370.1666
$ws.Cells.Item(16, 10).Value = 1900
$ws.Cells.Item(16, 11).Value = 1370.1666
$ws.Cells.Item(16, 12).Value = 1900
$ws.Cells.Item(16, 13).Value = -1083.1666
$ws.Cells.Item(16, 14).Value = -2474

# CRP row 99
$ws.Cells.Item(99, 8).Value = 1093.6666
$ws.Cells.Item(99, 9).Value = 1118.3334
$ws.Cells.Item(99, 10).Value = 1069
$ws.Cells.Item(99, 11).Value = 1118.3334
$ws.Cells.Item(99, 12).Value = 1069
$ws.Cells.Item(99, 13).Value = 379.6666
$ws.Cells.Item(99, 14).Value = -4065

# CRP row 105
$ws.Cells.Item(105, 8).Value = 938.8889
$ws.Cells.Item(105, 9).Value = 939.2308
$ws.Cells.Item(105, 10).Value = 938
$ws.Cells.Item(105, 11).Value = 939.2308
$ws.Cells.Item(105, 12).Value = 938
$ws.Cells.Item(105, 13).Value = 807.7692
$ws.Cells.Item(105, 14).Value = -4432

# CRP row 107
$ws.Cells.Item(107, 8).Value = 1875.1333
$ws.Cells.Item(107, 9).Value = 2801.5557
$ws.Cells.Item(107, 10).Value = 485.5
$ws.Cells.Item(107, 11).Value = 2801.5557
$ws.Cells.Item(107, 12).Value = 485.5
$ws.Cells.Item(107, 13).Value = -881.5556999999999
$ws.Cells.Item(107, 14).Value = -4325.5

# CRP row 113
$ws.Cells.Item(113, 8).Value = 1502.625
$ws.Cells.Item(113, 9).Value = 1370.1666
$ws.Cells.Item(113, 10).Value = 1900
$ws.Cells.Item(113, 11).Value = 1370.1666
$ws.Cells.Item(113, 12).Value = 1900
$ws.Cells.Item(113, 13).Value = 799.8334
$ws.Cells.Item(113, 14).Value = -6240

# CRP row 126
$ws.Cells.Item(126, 8).Value = 1093.6666
$ws.Cells.Item(126, 9).Value = 1118.3334
$ws.Cells.Item(126, 10).Value = 1069
$ws.Cells.Item(126, 11).Value = 3355.0002
$ws.Cells.Item(126, 12).Value = 3207
$ws.Cells.Item(126, 13).Value = -885.0001999999999
$ws.Cells.Item(126, 14).Value = -8147

$ws = $wb.Worksheets.Item("CUL")
# CUL row 110
$ws.Cells.Item(110, 8).Value = 2581.6667
$ws.Cells.Item(110, 9).Value = 2581.6667
$ws.Cells.Item(110, 11).Value = 7745.000100000001
$ws.Cells.Item(110, 13).Value = -3655.000100000001

# CUL row 113
$ws.Cells.Item(113, 8).Value = 3180.5789
$ws.Cells.Item(113, 9).Value = 549.5263
$ws.Cells.Item(113, 10).Value = 5811.6313
$ws.Cells.Item(113, 11).Value = 1648.5789
$ws.Cells.Item(113, 12).Value = 17434.8939
$ws.Cells.Item(113, 13).Value = 521.4211
$ws.Cells.Item(113, 14).Value = -21774.8939

# CUL row 133
$ws.Cells.Item(133, 8).Value = 6237.8237
$ws.Cells.Item(133, 10).Value = 7291.0835
$ws.Cells.Item(133, 12).Value = 21873.2505
$ws.Cells.Item(133, 14).Value = -31993.2505

$ws = $wb.Worksheets.Item("GSM")
# GSM row 97
$ws.Cells.Item(97, 8).Value = 799.6415
$ws.Cells.Item(97, 9).Value = 798.3570999999999
$ws.Cells.Item(97, 10).Value = 804.5454999999999
$ws.Cells.Item(97, 11).Value = 798.3570999999999
$ws.Cells.Item(97, 12).Value = 804.5454999999999
$ws.Cells.Item(97, 13).Value = -302.3570999999999
$ws.Cells.Item(97, 14).Value = -1796.5455

# GSM row 126
$ws.Cells.Item(126, 8).Value = 1581.2858
$ws.Cells.Item(126, 9).Value = 1793.3334
$ws.Cells.Item(126, 10).Value = 1422.25
$ws.Cells.Item(126, 11).Value = 5380.0002
$ws.Cells.Item(126, 12).Value = 4266.75
$ws.Cells.Item(126, 13).Value = -2910.0002
$ws.Cells.Item(126, 14).Value = -9206.75

$ws = $wb.Worksheets.Item("LTW")
# LTW row 7
$ws.Cells.Item(7, 8).Value = 2104.1
$ws.Cells.Item(7, 9).Value = 2289.2
$ws.Cells.Item(7, 10).Value = 1919
$ws.Cells.Item(7, 11).Value = 2289.2
$ws.Cells.Item(7, 12).Value = 1919
$ws.Cells.Item(7, 13).Value = -2177.2
$ws.Cells.Item(7, 14).Value = -2143

# LTW row 40
$ws.Cells.Item(40, 8).Value = 1486.7142
$ws.Cells.Item(40, 9).Value = 1429.125
$ws.Cells.Item(40, 10).Value = 1671
$ws.Cells.Item(40, 11).Value = 1429.125
$ws.Cells.Item(40, 12).Value = 1671
$ws.Cells.Item(40, 13).Value = -1293.125
$ws.Cells.Item(40, 14).Value = -1943

# LTW row 61
$ws.Cells.Item(61, 8).Value = 1282.138
$ws.Cells.Item(61, 9).Value = 1016.0476
$ws.Cells.Item(61, 10).Value = 1980.625
$ws.Cells.Item(61, 11).Value = 1016.0476
$ws.Cells.Item(61, 12).Value = 1980.625
$ws.Cells.Item(61, 13).Value = -814.0476
$ws.Cells.Item(61, 14).Value = -2384.625

# LTW row 113
$ws.Cells.Item(113, 8).Value = 1282.138
$ws.Cells.Item(113, 9).Value = 1016.0476
$ws.Cells.Item(113, 10).Value = 1980.625
$ws.Cells.Item(113, 11).Value = 1016.0476
$ws.Cells.Item(113, 12).Value = 1980.625
$ws.Cells.Item(113, 13).Value = 1153.9524
$ws.Cells.Item(113, 14).Value = -6320.625

# LTW row 122
$ws.Cells.Item(122, 8).Value = 6964.625
$ws.Cells.Item(122, 9).Value = 9235.143
$ws.Cells.Item(122, 10).Value = 3785.9
$ws.Cells.Item(122, 11).Value = 27705.429
$ws.Cells.Item(122, 12).Value = 11357.7
$ws.Cells.Item(122, 13).Value = -25255.429
$ws.Cells.Item(122, 14).Value = -16257.7

# LTW row 126
$ws.Cells.Item(126, 8).Value = 2104.1
$ws.Cells.Item(126, 9).Value = 2289.2
$ws.Cells.Item(126, 10).Value = 1919
$ws.Cells.Item(126, 11).Value = 6867.599999999999
$ws.Cells.Item(126, 12).Value = 5757
$ws.Cells.Item(126, 13).Value = -4397.599999999999
$ws.Cells.Item(126, 14).Value = -10697

# LTW row 132
$ws.Cells.Item(132, 8).Value = 6254.3105
$ws.Cells.Item(132, 9).Value = 8803.5625
$ws.Cells.Item(132, 10).Value = 3116.7693
$ws.Cells.Item(132, 11).Value = 26410.6875
$ws.Cells.Item(132, 12).Value = 9350.3079
$ws.Cells.Item(132, 13).Value = -23880.6875
$ws.Cells.Item(132, 14).Value = -14410.3079

# LTW row 136
$ws.Cells.Item(136, 8).Value = 2045.2307
$ws.Cells.Item(136, 9).Value = 987.25
$ws.Cells.Item(136, 10).Value = 3738
$ws.Cells.Item(136, 11).Value = 2961.75
$ws.Cells.Item(136, 12).Value = 11214
$ws.Cells.Item(136, 13).Value = -411.75
$ws.Cells.Item(136, 14).Value = -16314

$ws = $wb.Worksheets.Item("WVR")
# WVR row 96
$ws.Cells.Item(96, 8).Value = 16667642
$ws.Cells.Item(96, 9).Value = 50000950
$ws.Cells.Item(96, 10).Value = 987.5
$ws.Cells.Item(96, 11).Value = 50000950
$ws.Cells.Item(96, 12).Value = 987.5
$ws.Cells.Item(96, 13).Value = -49999577
$ws.Cells.Item(96, 14).Value = -3733.5
